# Apply attendance-log updates to Sheet1.
# - Corrects a few name/class mix-ups in the existing rows (2-13)
# - Adds new arrival/departure entries for a new student (Detlef Soost, 1a)
#   and additional check-in/check-out rows (14-20)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing rows (Vornamen / Nachnamen / Klasse corrections) ---

# Row 2: Eli Enders 4a -> Stephan Fuchs 3C
$ws.Range("A2").Value = "Stephan"
$ws.Range("B2").Value = "Fuchs"
$ws.Range("C2").Value = "3C"

# Row 3: Steve Mustermann 10x -> Eli Enders 4a
$ws.Range("A3").Value = "Eli"
$ws.Range("B3").Value = "Enders"
$ws.Range("C3").Value = "4a"

# Row 4: Stephan Fuchs 3C -> Max Schmitz 4a
$ws.Range("A4").Value = "Max"
$ws.Range("B4").Value = "Schmitz"
$ws.Range("C4").Value = "4a"

# Row 5: Steve Mustermann 10x -> Detlef Soost 1a
$ws.Range("A5").Value = "Detlef"
$ws.Range("B5").Value = "Soost"
$ws.Range("C5").Value = "1a"

# Row 6: Steve Mustermann 10x -> Stephan Fuchs 3C
$ws.Range("A6").Value = "Stephan"
$ws.Range("B6").Value = "Fuchs"
$ws.Range("C6").Value = "3C"

# Row 7 (Eli Enders 4a) is unchanged

# Row 8: Stephan Fuchs 3C -> Max Schmitz 4a
$ws.Range("A8").Value = "Max"
$ws.Range("B8").Value = "Schmitz"
$ws.Range("C8").Value = "4a"

# Rows 9 and 10 are unchanged

# Row 11: Steve Mustermann 10x -> Max Schmitz 4a
$ws.Range("A11").Value = "Max"
$ws.Range("B11").Value = "Schmitz"
$ws.Range("C11").Value = "4a"

# Row 12: Eli Enders 4a -> Max Schmitz (still 4a)
$ws.Range("A12").Value = "Max"
$ws.Range("B12").Value = "Schmitz"

# Row 13: Stephan Fuchs 3C -> Eli Enders 4a
$ws.Range("A13").Value = "Eli"
$ws.Range("B13").Value = "Enders"
$ws.Range("C13").Value = "4a"

# --- Append new log entries (rows 14-20) ---

$ws.Range("A14").Value = "Stephan"
$ws.Range("B14").Value = "Fuchs"
$ws.Range("C14").Value = "3C"
$ws.Range("D14").Value = "10.14.2023 20:28"
$ws.Range("E14").Value = "10.14.2023 20:28"

$ws.Range("A15").Value = "Detlef"
$ws.Range("B15").Value = "Soost"
$ws.Range("C15").Value = "1a"
$ws.Range("D15").Value = "10.14.2023 20:28"
$ws.Range("E15").Value = "10.14.2023 20:28"

$ws.Range("A16").Value = "Max"
$ws.Range("B16").Value = "Schmitz"
$ws.Range("C16").Value = "4a"
$ws.Range("D16").Value = "10.14.2023 20:59"
$ws.Range("E16").Value = "10.14.2023 20:59"

$ws.Range("A17").Value = "Detlef"
$ws.Range("B17").Value = "Soost"
$ws.Range("C17").Value = "1a"
$ws.Range("D17").Value = "10.15.2023 18:18"
$ws.Range("E17").Value = "10.15.2023 20:05"

$ws.Range("A18").Value = "Detlef"
$ws.Range("B18").Value = "Soost"
$ws.Range("C18").Value = "1a"
$ws.Range("D18").Value = "10.17.2023 18:56"
$ws.Range("E18").Value = "10.17.2023 19:57"

$ws.Range("A19").Value = "Max"
$ws.Range("B19").Value = "Schmitz"
$ws.Range("C19").Value = "4a"
$ws.Range("D19").Value = "10.17.2023 18:58"
$ws.Range("E19").Value = "10.17.2023 19:57"

$ws.Range("A20").Value = "Stephan"
$ws.Range("B20").Value = "Fuchs"
$ws.Range("C20").Value = "3C"
$ws.Range("D20").Value = "10.17.2023 19:16"
$ws.Range("E20").Value = "10.17.2023 19:57"
